$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.435.52"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.937.59"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.06"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.30"
$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -0.40%  "

$ws.Range("E9").Value = "  +2.63%  "

$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.71"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.467.74"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.383.52"

$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.934.02"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "434.49"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.53"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("E21").Value = "  -0.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.18"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.00"
$ws.Range("E23").Value = "  +1.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.98"
$ws.Range("E24").Value = "  +1.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  -2.46%  "

$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  -4.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.62"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("E31").Value = "  +3.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.80"
$ws.Range("E32").Value = "  +0.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0888"
$ws.Range("E34").Value = "  +2.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.02"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.03"
$ws.Range("E37").Value = "  -2.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.04"
$ws.Range("E38").Value = "  +0.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.123"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("E41").Value = "  +6.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.285"
$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0350"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "372.88"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.716.75"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.51"
$ws.Range("E46").Value = "  +2.84%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.02"
$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.02"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("E51").Value = "  -0.09%  "
